# Attendance up to date : Oct 11, 2021
#
# The SPONSOR sheet's "Oct" week column (F) had no data entered yet for the
# attendees. This edit fills in the new attendance week (column F, date 11)
# for all sponsor meeting attendees, mirroring the month/day header already
# present in column E, and marks the SPONSOR sheet/cell F11 as the
# last-active selection (reflecting where the user was working).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPONSOR")

$checkmark = [char]0x2714

# Header rows: month + day for the new attendance column (F), matching
# the pattern already used in column E ("Oct").
$ws.Range("F3").Value = "Oct"
$ws.Range("F4").Value = 11

# Mark attendance (checkmark) for every attendee row that meets this week.
$ws.Range("F5").Value = $checkmark
$ws.Range("F6").Value = $checkmark
$ws.Range("F7").Value = $checkmark
$ws.Range("F8").Value = $checkmark
$ws.Range("F9").Value = $checkmark
$ws.Range("F10").Value = $checkmark
$ws.Range("F11").Value = $checkmark

# The checkbox/attendance validation that used to cover F5:N11 as part of
# the "C5:D11 F5:N11" rule now excludes column F (it moves into the
# E5:F11 rule below), so split it off column F first.
$ws.Range("F5:F11").Validation.Delete()

# Column F (rows 5:11) now shares the same checkbox-list validation as
# column E, so fold it into that rule's range (E5:E11 -> E5:F11).
$rngEF = $ws.Range("E5:F11")
$rngEF.Validation.Delete()
$rngEF.Validation.Add(3, 1, 1, "TEAM!checkbox", "0")
$rngEF.Validation.IgnoreBlank = $true
$rngEF.Validation.InCellDropdown = $true
$rngEF.Validation.ShowInput = $true
$rngEF.Validation.ShowError = $true

# SPONSOR becomes the active sheet/tab, with F11 as the selected cell.
$ws.Activate() | Out-Null
$ws.Range("F11").Select() | Out-Null
